# Apply the "Questions + Tasks Update + Bucket's Teeth" commit changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Tasks update / Bucket's Teeth: completion percentages for Task 3 (Bucket) ---
# B28 = "Attachments" -> completion 100 -> 50
# B29 = "Teeths"       -> completion 0   -> 100
$ws.Range("D28").Value = 50
$ws.Range("D29").Value = 100

# --- Rename "Boom Attachment" to "Boom's Attachment" (row 45, Task 7 / Arm&Boom) ---
$ws.Range("B45").Value = "Boom's Attachment"

# --- Update the saved view/selection to match the author's last position ---
$ws.Activate()
$ws.Range("D30").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
